# Update TestResult.xlsx "Sheet1":
#  - Row 4 (Click on Apply Now) now PASSes, so its failure Remarks (D4) is cleared.
#  - Rows 7-14 keep their FAIL result, but the Selenium "no such element" messages
#    are replaced with a "chrome not reachable" message (the browser session died).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4: Click on Apply Now -> PASS, remove the old error remark entirely.
$ws.Range("C4").Value = "PASS"
$ws.Range("D4").ClearContents()

# New shared message for rows whose element lookups now fail earlier,
# because Chrome itself became unreachable.
$chromeNotReachable = "Message: chrome not reachable`n  (Session info: chrome=91.0.4472.114)`n"

$ws.Range("D7").Value = $chromeNotReachable
$ws.Range("D8").Value = $chromeNotReachable
$ws.Range("D9").Value = $chromeNotReachable
$ws.Range("D10").Value = $chromeNotReachable
$ws.Range("D11").Value = $chromeNotReachable
$ws.Range("D12").Value = $chromeNotReachable
$ws.Range("D13").Value = $chromeNotReachable
$ws.Range("D14").Value = $chromeNotReachable
